$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-14 Friday" "2025-03-15 Saturday"

Replace-Text "640÷8=80, 0" "999÷5=199, 4"
Replace-Text "656÷9=72, 8" "672÷2=336, 0"
Replace-Text "771÷2=385, 1" "433÷2=216, 1"
Replace-Text "996÷2=498, 0" "670÷7=95, 5"
Replace-Text "352÷2=176, 0" "229÷7=32, 5"

Replace-Text "128÷8=16, 0" "655÷3=218, 1"
Replace-Text "113÷3=37, 2" "730÷4=182, 2"
Replace-Text "506÷3=168, 2" "176÷5=35, 1"
Replace-Text "537÷4=134, 1" "690÷5=138, 0"
Replace-Text "290÷4=72, 2" "834÷7=119, 1"

Replace-Text "455÷6=75, 5" "539÷2=269, 1"
Replace-Text "816÷5=163, 1" "956÷6=159, 2"
Replace-Text "363÷2=181, 1" "940÷9=104, 4"
Replace-Text "304÷5=60, 4" "480÷8=60, 0"
Replace-Text "967÷3=322, 1" "738÷4=184, 2"

Replace-Text "864÷4=216, 0" "591÷7=84, 3"
Replace-Text "143÷9=15, 8" "673÷6=112, 1"
Replace-Text "613÷6=102, 1" "229÷2=114, 1"
Replace-Text "537÷3=179, 0" "741÷9=82, 3"
Replace-Text "695÷9=77, 2" "301÷8=37, 5"

Replace-Text "342÷9=38, 0" "911÷7=130, 1"
Replace-Text "119÷9=13, 2" "265÷2=132, 1"
Replace-Text "434÷8=54, 2" "297÷2=148, 1"
Replace-Text "645÷5=129, 0" "533÷2=266, 1"
Replace-Text "695÷3=231, 2" "674÷9=74, 8"
